$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 corresponds to participant 5 (session data collection)
# C6: Year since course
$ws.Range("C6").Value = 1

# E6: date of session 1
$ws.Range("E6").Value = 20250109

# H6: comments
$ws.Range("H6").Value = "TA for the lab course (2023, 2024)"

# Update the active cell selection to F6 as recorded in the session
$ws.Range("F6").Select()
